# Add a new "posts" worksheet (sheetId=3) after the existing "creators" sheet,
# populate it with post data, apply header-bold + a thin-left border accent on
# the second row, and move the active tab / selections to match.

$wb = $excel.ActiveWorkbook

# --- Insert the new sheet at the end of the tab strip ---------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$posts = $wb.Worksheets.Add($null, $lastSheet)
$posts.Name = "posts"

# --- Header row -------------------------------------------------------------
$posts.Range("A1").Value = "title"
$posts.Range("B1").Value = "content"
$posts.Range("C1").Value = "author"
$posts.Range("D1").Value = "image"
$posts.Range("A1:D1").Font.Bold = $true

# --- Row 2: The Boat ---------------------------------------------------------
$posts.Range("A2").Value = "The Boat"
$posts.Range("B2").Value = "Look at this boat!"
$posts.Range("C2").Value = "BobbyPaints"
$posts.Range("D2").Value = "sailing_boat.jpg"
$posts.Range("C2:D2").Borders.Item(11).LineStyle = 1

# --- Row 3: Dolphin -----------------------------------------------------------
$posts.Range("A3").Value = "Dolphin"
$posts.Range("B3").Value = "Cubism or something!"
$posts.Range("C3").Value = "BobbyPaints"
$posts.Range("D3").Value = "dolphin.jpg"

# --- Row 4: Year of the Monkey (no content/author text cell) ----------------
$posts.Range("A4").Value = "Year of the Monkey"
$posts.Range("C4").Value = "JohnnyDraws"
$posts.Range("D4").Value = "monkey.jpg"

# --- Row 5: My House ----------------------------------------------------------
$posts.Range("A5").Value = "My House"
$posts.Range("B5").Value = "I can paint too!"
$posts.Range("C5").Value = "HelenSculpts"
$posts.Range("D5").Value = "house.jpg"

# --- Row 6: Owls are great! ---------------------------------------------------
$posts.Range("A6").Value = "Owls are great!"
$posts.Range("B6").Value = "superb, even!"
$posts.Range("C6").Value = "JohnnyDraws"
$posts.Range("D6").Value = "owl.jpg"

# --- Selections / active tab -------------------------------------------------
$creators = $wb.Worksheets.Item("creators")
$creators.Activate() | Out-Null
$creators.Range("A5").Select() | Out-Null

$posts.Activate() | Out-Null
$posts.Range("D6").Select() | Out-Null
